$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 65, pushing existing rows 65..174 down to 67..176.
$ws.Rows("65:66").Insert()

# Populate the two newly-inserted rows with their data (same market/product metadata
# as the rest of the sheet, new Fecha/Calidad/Volumen/Precio values).

# Row 65: Maracuyá, Especial -> became Primera per new record
$ws.Range("A65").Value = 1
$ws.Range("B65").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C65").Value = "Arica y Parinacota"
$ws.Range("D65").Value = 45070
$ws.Range("E65").Value = 15
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100108
$ws.Range("H65").Value = "Tropicales y subtropicales"
$ws.Range("I65").Value = 100108003
$ws.Range("J65").Value = "Maracuyá"
$ws.Range("K65").Value = "Sin especificar"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 120
$ws.Range("N65").Value = 34000
$ws.Range("O65").Value = 35000
$ws.Range("P65").Value = 34500
$ws.Range("Q65").Value = "$/caja 20 kilos"
$ws.Range("R65").Value = "Región de Arica y Parinacota"
$ws.Range("S65").Value = 1725
$ws.Range("T65").Value = 20

# Row 66
$ws.Range("A66").Value = 1
$ws.Range("B66").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C66").Value = "Arica y Parinacota"
$ws.Range("D66").Value = 45070
$ws.Range("E66").Value = 15
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100108
$ws.Range("H66").Value = "Tropicales y subtropicales"
$ws.Range("I66").Value = 100108003
$ws.Range("J66").Value = "Maracuyá"
$ws.Range("K66").Value = "Sin especificar"
$ws.Range("L66").Value = "Segunda"
$ws.Range("M66").Value = 140
$ws.Range("N66").Value = 27000
$ws.Range("O66").Value = 28000
$ws.Range("P66").Value = 27500
$ws.Range("Q66").Value = "$/caja 20 kilos"
$ws.Range("R66").Value = "Región de Arica y Parinacota"
$ws.Range("S66").Value = 1375
$ws.Range("T66").Value = 20

# Ensure the date cells keep the same number format as the rest of the Fecha column.
$ws.Range("D65:D66").NumberFormat = $ws.Range("D67").NumberFormat
